$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPEC CPU2006")

# Row 6 - bwaves_r
$ws.Range("C6").Value = 1753042400
$ws.Range("D6").Value = 1747289800
$ws.Range("F6").Value = 0.16
$ws.Range("G6").Value = 0.16
$ws.Range("I6").Value = 28221816
$ws.Range("J6").Value = 28219798
$ws.Range("L6").Value = 9719732
$ws.Range("M6").Value = 9719861
$ws.Range("O6").Value = 18502084
$ws.Range("P6").Value = 9719861
$ws.Range("R6").Value = 289095
$ws.Range("S6").Value = 151872
$ws.Range("U6").Value = 850
$ws.Range("V6").Value = 848

# Row 24 - roms_r
$ws.Range("C24").Value = 836215000
$ws.Range("D24").Value = 836629600
$ws.Range("F24").Value = 0.22
$ws.Range("G24").Value = 0.22
$ws.Range("I24").Value = 18999065
$ws.Range("J24").Value = 19002462
$ws.Range("L24").Value = 9238432
$ws.Range("M24").Value = 9242938
$ws.Range("O24").Value = 9760633
$ws.Range("P24").Value = 9242938
$ws.Range("R24").Value = 152509
$ws.Range("S24").Value = 144420
$ws.Range("U24").Value = 1162
$ws.Range("V24").Value = 1163
